$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'288.37"
$ws.Range("D3").Value = "'29.22"
$ws.Range("E3").Value = "'2.19%"
$ws.Range("D4").Value = "'5.098"
$ws.Range("E4").Value = "'3.67%"
$ws.Range("E5").Value = "'2.97%"
$ws.Range("D6").Value = "'7.331"
$ws.Range("E6").Value = "'1.54%"
$ws.Range("D7").Value = "'3.406"
$ws.Range("E7").Value = "'1.16%"
$ws.Range("E8").Value = "'1.24%"
$ws.Range("D9").Value = "'0.9185"
$ws.Range("E9").Value = "'0.53%"
$ws.Range("D10").Value = "'0.1587"
$ws.Range("E10").Value = "'2.83%"
$ws.Range("D11").Value = "'0.06810"
$ws.Range("E11").Value = "'7.41%"
$ws.Range("D12").Value = "'0.07665"
$ws.Range("E12").Value = "'0.26%"
$ws.Range("D13").Value = "'0.02936"
$ws.Range("E13").Value = "'-1.34%"
$ws.Range("D14").Value = "'0.08985"
$ws.Range("E14").Value = "'0.22%"
$ws.Range("D15").Value = "'0.001572"
$ws.Range("E15").Value = "'-1.71%"
$ws.Range("E16").Value = "'0.84%"
$ws.Range("E17").Value = "'-0.80%"
$ws.Range("D18").Value = "'0.006253"
$ws.Range("E18").Value = "'3.03%"
$ws.Range("E19").Value = "'-0.30%"
$ws.Range("D20").Value = "'2.220"
$ws.Range("E21").Value = "'2.02%"
$ws.Range("E22").Value = "'-2.47%"
$ws.Range("D23").Value = "'4.076"
$ws.Range("E23").Value = "'1.98%"
$ws.Range("E25").Value = "'0.11%"
$ws.Range("D26").Value = "'0.004115"
$ws.Range("E26").Value = "'-4.94%"
$ws.Range("E27").Value = "'1.62%"
$ws.Range("D28").Value = "'0.0001617"
$ws.Range("E28").Value = "'-1.15%"
$ws.Range("D40").Value = "'0.04225"
$ws.Range("E40").Value = "'1.54%"
$ws.Range("D41").Value = "'0.006731"
$ws.Range("E41").Value = "'0.35%"
$ws.Range("E42").Value = "'0.47%"
$ws.Range("E44").Value = "'13.92%"
$ws.Range("D45").Value = "'0.00005708"
$ws.Range("E45").Value = "'6.30%"
$ws.Range("D46").Value = "'1.974"
$ws.Range("E46").Value = "'-3.28%"
$ws.Range("E47").Value = "'-29.41%"
